$wb = $excel.ActiveWorkbook

# Books sheet: reduce Quantity for a few items
$ws = $wb.Worksheets.Item("Books")
$ws.Range("D4").Value = 5   # In Search of Lost Time
$ws.Range("D8").Value = 7   # The Red and the Black
$ws.Range("D11").Value = 7  # Leaves of Grass

# Fruits sheet: reduce Quantity for a couple of items
$ws = $wb.Worksheets.Item("Fruits")
$ws.Range("D5").Value = 4   # Orange
$ws.Range("D6").Value = 8   # Gava

# Games sheet: reduce Quantity for a couple of items
$ws = $wb.Worksheets.Item("Games")
$ws.Range("D2").Value = 8   # Far Cry 6
$ws.Range("D11").Value = 8  # Rayman Legends

# Sports Kit sheet: reduce Quantity for a couple of items
$ws = $wb.Worksheets.Item("Sports Kit")
$ws.Range("D2").Value = 9   # Cricket
$ws.Range("D9").Value = 8   # Badminton

$wb.Save()
